$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PoiFormulaHelperTest")

# Update the formula in A9: change the second SUM range from A1:D5 to B1:D5
$ws.Range("A9").Formula = "=SUM(Data!A1:D5)+SUM(Data!B1:D5)"

# Update the selection shown in the sheet view to O39
$ws.Range("O39").Select()
